# Expand the "Conclusion" section of the lab report.
#
# The document currently ends with a single paragraph that holds nothing
# but a lone tab character plus the _GoBack bookmark (right before the
# section properties). We grow that into five fully-written conclusion
# paragraphs (each first-line indented like the rest of the write-up),
# keep the bookmark on the final one, and leave one bare trailing
# paragraph behind it, matching the reference document.

$d = $word.ActiveDocument

# Locate the placeholder paragraph: it's the last paragraph in the main
# story and it is the one that owns the _GoBack bookmark. It will become
# the paragraph holding the *last* piece of new text (so the bookmark
# stays attached to it), exactly like the reference document.
$target = $d.Paragraphs.Last
$targetIndex = $d.Paragraphs.Count

$texts = @(
    "Reference strings are a way to represent the unknown path of page requests a program will take along its execution.  Attempting to artificially generate locality taught be about the characteristic patterns that program execution usually follows. It is not common for a program to execution similar to the random locality reference strings.  Especially with modern compilers, code if highly optimized so that optimized replacement algorithms can optimize paging performance.",
    "Considering the replacement algorithms assigned, I realize why operating system programming today in age is impossible to be done with less than a few hundred or thousand people.  Every single decision involves a give and take.",
    "Coding a demand paging memory management system was very fun.  Then adding a round robin scheduler on top of that proved to be a very fun and stressful week.  The most important thing I learned is the dependence that an Operating System has on memory.  Due to the huge performance penalty, an Operating System must try to keep as much as possible in memory, but also be the most efficient about it as possible. Whatever memory the OS uses, it is memory it cannot use to increase performance of user programs.",
    "Operating System code execution is unavoidable.  However, by coding my own scheduling algorithm, I was able to catch an abstract glimpse of all the possible duties an OS must perform all the time in between user programs.  I wish I had more time to continue building upon this project. ",
    "I have three planned next steps.  The first is to code a small chance that a process will request a page outside its address space.  My simulation would have to catch that and terminate the program. Another idea is to create other types of resources other than memory that reference strings can request.  This would allow me to consider and code deadlock prevention and avoidance algorithms.  Finally, I would have like to have implemented a page table.  Due to time constraints, I implemented an inverse page table.  "
)

# Step 1: while the placeholder paragraph is still plain (no direct
# formatting yet), give it a sibling right after it - this will become
# the bare trailing paragraph and, being created now, inherits no
# first-line indent.
$target.Range.InsertParagraphAfter()

# Step 2: insert one new blank paragraph ahead of the placeholder for
# each of the first N-1 pieces of text; the placeholder itself becomes
# the paragraph for the last piece of text, so it keeps sliding down
# and lands right before the bare trailing paragraph from step 1.
for ($i = 0; $i -lt ($texts.Length - 1); $i++) {
    $target.Range.InsertParagraphBefore()
}

# $targetIndex was the bookmarked paragraph's slot before any inserts;
# the freshly-made blank paragraphs now occupy that slot and the ones
# right after it, with the (shifted-down) bookmarked paragraph landing
# right before the bare trailing paragraph.
$firstIndex = $targetIndex

for ($i = 0; $i -lt $texts.Length; $i++) {
    $p = $d.Paragraphs.Item($firstIndex + $i)
    $p.Range.ParagraphFormat.FirstLineIndent = 36
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Text = $texts[$i]
}

# The bare trailing paragraph from step 1 is the next one; round-trip a
# throwaway character through it so the engine materializes it cleanly
# (it is already un-indented and still otherwise untouched).
$trailingIndex = $firstIndex + $texts.Length
$trailing = $d.Paragraphs.Item($trailingIndex)
$tr = $d.Range($trailing.Range.Start, $trailing.Range.End - 1)
$tr.Text = "X"
$trailing2 = $d.Paragraphs.Item($trailingIndex)
$tr2 = $d.Range($trailing2.Range.Start, $trailing2.Range.End - 1)
$tr2.Delete()

Write-Host "done"
